$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows right before the current row 392, pushing the existing
# rows 392-407 down to 396-411 (their contents, and relative formatting,
# move with them automatically).
$ws.Range("A392:R395").EntireRow.Insert()

# New weekly records for fecha 44568 (2022-01-07), inserted at the top of
# this Repollo price block.
$newRows = @(
    @{ Row=392; H="Copenhague";     I="Primera"; J=1600; K=700; L=800; M=750; O="Región Metropolitana" },
    @{ Row=393; H="Crespo record";  I="Primera"; J=3400; K=800; L=900; M=850; O="Región Metropolitana" },
    @{ Row=394; H="Crespo record";  I="Segunda"; J=1600; K=600; L=600; M=600; O="Región Metropolitana" },
    @{ Row=395; H="Morada(o)";      I="Primera"; J=1600; K=800; L=900; M=850; O="Región Metropolitana" }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value2 = 9
    $ws.Cells.Item($row, 2).Value2 = "Vega Central Mapocho de Santiago"
    $ws.Cells.Item($row, 3).Value2 = "Metropolitana"
    $ws.Cells.Item($row, 4).Value2 = 44568
    $ws.Cells.Item($row, 5).Value2 = 13
    $ws.Cells.Item($row, 6).Value2 = 100112006
    $ws.Cells.Item($row, 7).Value2 = "Repollo"
    $ws.Cells.Item($row, 8).Value2 = $r.H
    $ws.Cells.Item($row, 9).Value2 = $r.I
    $ws.Cells.Item($row, 10).Value2 = $r.J
    $ws.Cells.Item($row, 11).Value2 = $r.K
    $ws.Cells.Item($row, 12).Value2 = $r.L
    $ws.Cells.Item($row, 13).Value2 = $r.M
    $ws.Cells.Item($row, 14).Value2 = "`$/unidad"
    $ws.Cells.Item($row, 15).Value2 = $r.O
    $ws.Cells.Item($row, 16).Value2 = $r.M
    $ws.Cells.Item($row, 17).Value2 = 1
    $ws.Cells.Item($row, 18).Value2 = "Hortaliza"
}
